$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45204) for every data
# row (2..536). The update bumps that value by one day (45204 -> 45205)
# across the whole column range.
$ws.Range("C2:C536").Value = 45205
